$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.339.07"
$ws.Range("E2").Value = "  +0.45%  "

# Row 3
$ws.Range("D3").Value = "1.788.15"
$ws.Range("E3").Value = "  +1.72%  "

# Row 4
$ws.Range("E4").Value = "  +0.95%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.39%  "

# Row 6
$ws.Range("E6").Value = "  +0.88%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3799"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.34%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3449"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.85%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.68"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.37%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.197"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.31%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07490"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.17%  "

# Row 12
$ws.Range("E12").Value = "  +0.91%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.20%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.458"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.57%  "

# Row 15
$ws.Range("D15").Value = "1.790.76"
$ws.Range("E15").Value = "  +2.02%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.069"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.44%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001102"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.58%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06680"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.69%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "84.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.48%  "

# Row 20
$ws.Range("E20").Value = "  +0.72%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.538"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.47%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.29%  "

# Row 23
$ws.Range("D23").Value = "27.352.00"
$ws.Range("E23").Value = "  +0.63%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.34%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.430"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.39%  "

# Row 26
$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.502"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.54%  "

# Row 27
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.557"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.40%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.76%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "153.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.63%  "

# Row 30
$ws.Range("D30").Value = "1.993.77"
$ws.Range("E30").Value = "  +2.04%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "133.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.31%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.060"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.09%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.071"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.43%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08686"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.55%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.44%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.649"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.49%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.475"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.51%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6905"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.41%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06389"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.66%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2202"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.17%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.851"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.92%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02342"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.31%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.260"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.63%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.02%  "

# Row 45
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6456"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.67%  "

# Row 46
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.002"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.69%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.862"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.07%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.131"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.46%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "129.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.67%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07192"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.85%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.69%  "
